$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 618.84
$ws.Range("I28").Value = 539.3889
$ws.Range("J28").Value = 823.1429000000001
$ws.Range("K28").Value = 539.3889
$ws.Range("L28").Value = 823.1429000000001
$ws.Range("M28").Value = -54.38890000000004
$ws.Range("N28").Value = -1793.1429

$ws.Range("H113").Value = 3655.0557
$ws.Range("I113").Value = 3290.3333
$ws.Range("J113").Value = 4019.7778
$ws.Range("K113").Value = 3290.3333
$ws.Range("L113").Value = 4019.7778
$ws.Range("M113").Value = -36.33329999999978
$ws.Range("N113").Value = -10527.7778

$ws.Range("H115").Value = 2794.3076
$ws.Range("I115").Value = 632.6
$ws.Range("J115").Value = 10000
$ws.Range("K115").Value = 1897.8
$ws.Range("L115").Value = 30000
$ws.Range("M115").Value = -330.8000000000002
$ws.Range("N115").Value = -33134

$ws.Range("H132").Value = 1068346
$ws.Range("I132").Value = 3047.7
$ws.Range("K132").Value = 9143.099999999999
$ws.Range("M132").Value = -6613.099999999999

$ws.Range("H137").Value = 1820613.4
$ws.Range("I137").Value = 2942796.5
$ws.Range("J137").Value = 3745.3809
$ws.Range("K137").Value = 8828389.5
$ws.Range("L137").Value = 11236.1427
$ws.Range("M137").Value = -8825839.5
$ws.Range("N137").Value = -16336.1427

$ws.Range("H138").Value = 3281526.2
$ws.Range("I138").Value = 3209.818
$ws.Range("J138").Value = 4002755.8
$ws.Range("K138").Value = 9629.454000000002
$ws.Range("L138").Value = 12008267.4
$ws.Range("M138").Value = -4489.454000000002
$ws.Range("N138").Value = -12018547.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 3013
$ws.Range("I4").Value = 200
$ws.Range("J4").Value = 3950.6667
$ws.Range("K4").Value = 200
$ws.Range("L4").Value = 3950.6667
$ws.Range("M4").Value = -84
$ws.Range("N4").Value = -4182.6667

$ws.Range("H45").Value = 1496.5
$ws.Range("I45").Value = 1444.4
$ws.Range("J45").Value = 1583.3334
$ws.Range("K45").Value = 1444.4
$ws.Range("L45").Value = 1583.3334
$ws.Range("M45").Value = -1067.4
$ws.Range("N45").Value = -2337.3334

$ws.Range("H74").Value = 7137190.5
$ws.Range("I74").Value = 13375165
$ws.Range("J74").Value = 48583.453
$ws.Range("K74").Value = 13375165
$ws.Range("L74").Value = 48583.453
$ws.Range("M74").Value = -13374291
$ws.Range("N74").Value = -50331.453

$ws.Range("H77").Value = 7137190.5
$ws.Range("I77").Value = 13375165
$ws.Range("J77").Value = 48583.453
$ws.Range("K77").Value = 66875825
$ws.Range("L77").Value = 242917.265
$ws.Range("M77").Value = -66871457
$ws.Range("N77").Value = -251653.265

$ws.Range("H118").Value = 59875
$ws.Range("J118").Value = 59875
$ws.Range("L118").Value = 59875
$ws.Range("N118").Value = -63189

$ws.Range("H122").Value = 13891701
$ws.Range("I122").Value = 2262
$ws.Range("J122").Value = 27781140
$ws.Range("K122").Value = 6786
$ws.Range("L122").Value = 83343420
$ws.Range("M122").Value = -4336
$ws.Range("N122").Value = -83348320

$ws.Range("H132").Value = 49292.58
$ws.Range("I132").Value = 32213.734
$ws.Range("J132").Value = 102081.73
$ws.Range("K132").Value = 96641.202
$ws.Range("L132").Value = 306245.19
$ws.Range("M132").Value = -94111.202
$ws.Range("N132").Value = -311305.19

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1772.4445
$ws.Range("I20").Value = 1657.3334
$ws.Range("J20").Value = 2002.6666
$ws.Range("K20").Value = 1657.3334
$ws.Range("L20").Value = 2002.6666
$ws.Range("M20").Value = -1410.3334
$ws.Range("N20").Value = -2496.6666

$ws.Range("H94").Value = 472.5
$ws.Range("I94").Value = 425.9375
$ws.Range("J94").Value = 845
$ws.Range("K94").Value = 425.9375
$ws.Range("L94").Value = 845
$ws.Range("M94").Value = 25.0625
$ws.Range("N94").Value = -1747

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1602946
$ws.Range("I22").Value = 2315075
$ws.Range("J22").Value = 655.5
$ws.Range("K22").Value = 2315075
$ws.Range("L22").Value = 655.5
$ws.Range("M22").Value = -2314725
$ws.Range("N22").Value = -1355.5

$ws.Range("H31").Value = 248954.66
$ws.Range("I31").Value = 54475.527
$ws.Range("J31").Value = 372124.78
$ws.Range("K31").Value = 54475.527
$ws.Range("L31").Value = 372124.78
$ws.Range("M31").Value = -54180.527
$ws.Range("N31").Value = -372714.78

$ws.Range("H34").Value = 248954.66
$ws.Range("I34").Value = 54475.527
$ws.Range("J34").Value = 372124.78
$ws.Range("K34").Value = 54475.527
$ws.Range("L34").Value = 372124.78
$ws.Range("M34").Value = -54273.527
$ws.Range("N34").Value = -372528.78

$ws.Range("H52").Value = 50000
$ws.Range("J52").Value = 50000
$ws.Range("L52").Value = 50000
$ws.Range("N52").Value = -50588

$ws.Range("H107").Value = 524.7857
$ws.Range("I107").Value = 466.26923
$ws.Range("J107").Value = 619.875
$ws.Range("K107").Value = 466.26923
$ws.Range("L107").Value = 619.875
$ws.Range("M107").Value = 1453.73077
$ws.Range("N107").Value = -4459.875

$ws.Range("H122").Value = 2314.7727
$ws.Range("I122").Value = 1918.3334
$ws.Range("J122").Value = 2790.5
$ws.Range("K122").Value = 5755.0002
$ws.Range("L122").Value = 8371.5
$ws.Range("M122").Value = -3305.0002
$ws.Range("N122").Value = -13271.5

$ws.Range("H129").Value = 37833
$ws.Range("J129").Value = 37833
$ws.Range("L129").Value = 37833
$ws.Range("N129").Value = -47833

$ws.Range("H130").Value = 67295
$ws.Range("J130").Value = 67295
$ws.Range("L130").Value = 67295
$ws.Range("N130").Value = -77335

$ws.Range("H131").Value = 35000
$ws.Range("J131").Value = 35000
$ws.Range("L131").Value = 35000
$ws.Range("N131").Value = -45080

$ws.Range("H132").Value = 93864.45
$ws.Range("I132").Value = 2418.5
$ws.Range("J132").Value = 203599.6
$ws.Range("K132").Value = 7255.5
$ws.Range("L132").Value = 610798.8
$ws.Range("M132").Value = -4725.5
$ws.Range("N132").Value = -615858.8

$ws.Range("H134").Value = 55089.9
$ws.Range("I134").Value = 1000
$ws.Range("K134").Value = 3000
$ws.Range("M134").Value = -465

$ws.Range("H141").Value = 47693.777
$ws.Range("J141").Value = 47693.777
$ws.Range("L141").Value = 47693.777
$ws.Range("N141").Value = -58053.777

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 2428.3125
$ws.Range("I109").Value = 975.5
$ws.Range("K109").Value = 2926.5
$ws.Range("M109").Value = -1886.5

$ws.Range("H113").Value = 626.1852
$ws.Range("I113").Value = 561.8421
$ws.Range("J113").Value = 779
$ws.Range("K113").Value = 1685.5263
$ws.Range("L113").Value = 2337
$ws.Range("M113").Value = 484.4737
$ws.Range("N113").Value = -6677

$ws.Range("H122").Value = 1009.15625
$ws.Range("I122").Value = 310.6
$ws.Range("J122").Value = 1138.5186
$ws.Range("K122").Value = 2795.4
$ws.Range("L122").Value = 10246.6674
$ws.Range("M122").Value = -345.4000000000001
$ws.Range("N122").Value = -15146.6674

$ws.Range("H130").Value = 2857.1177
$ws.Range("I130").Value = 913.6667
$ws.Range("J130").Value = 3273.5715
$ws.Range("K130").Value = 2741.0001
$ws.Range("L130").Value = 9820.7145
$ws.Range("M130").Value = 2278.9999
$ws.Range("N130").Value = -19860.7145

$ws.Range("H131").Value = 855
$ws.Range("J131").Value = 1029.8064
$ws.Range("L131").Value = 3089.4192
$ws.Range("N131").Value = -13169.4192

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2894
$ws.Range("I102").Value = 3090.8572
$ws.Range("J102").Value = 2434.6667
$ws.Range("K102").Value = 3090.8572
$ws.Range("L102").Value = 2434.6667
$ws.Range("M102").Value = -1468.8572
$ws.Range("N102").Value = -5678.6667

$ws.Range("H122").Value = 2446.3076
$ws.Range("I122").Value = 1977.1111
$ws.Range("J122").Value = 3502
$ws.Range("K122").Value = 5931.3333
$ws.Range("L122").Value = 10506
$ws.Range("M122").Value = -3481.3333
$ws.Range("N122").Value = -15406

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 500.85715
$ws.Range("I22").Value = 409.13635
$ws.Range("J22").Value = 601.75
$ws.Range("K22").Value = 409.13635
$ws.Range("L22").Value = 601.75
$ws.Range("M22").Value = -114.13635
$ws.Range("N22").Value = -1191.75

$ws.Range("H27").Value = 500.85715
$ws.Range("I27").Value = 409.13635
$ws.Range("J27").Value = 601.75
$ws.Range("K27").Value = 409.13635
$ws.Range("L27").Value = 601.75
$ws.Range("M27").Value = -302.13635
$ws.Range("N27").Value = -815.75

$ws.Range("H82").Value = 2188.7
$ws.Range("I82").Value = 1498.75
$ws.Range("J82").Value = 2648.6667
$ws.Range("K82").Value = 1498.75
$ws.Range("L82").Value = 2648.6667
$ws.Range("M82").Value = -1137.75
$ws.Range("N82").Value = -3370.6667

$ws.Range("H85").Value = 2188.7
$ws.Range("I85").Value = 1498.75
$ws.Range("J85").Value = 2648.6667
$ws.Range("K85").Value = 1498.75
$ws.Range("L85").Value = 2648.6667
$ws.Range("M85").Value = -250.75
$ws.Range("N85").Value = -5144.6667

$ws.Range("H100").Value = 1964.7778
$ws.Range("I100").Value = 1920.75
$ws.Range("K100").Value = 1920.75
$ws.Range("M100").Value = -1379.75

$ws.Range("H122").Value = 3546.2683
$ws.Range("I122").Value = 3997.4707
$ws.Range("J122").Value = 3226.6667
$ws.Range("K122").Value = 11992.4121
$ws.Range("L122").Value = 9680.000100000001
$ws.Range("M122").Value = -9542.4121
$ws.Range("N122").Value = -14580.0001

$ws.Range("H132").Value = 80715.12
$ws.Range("I132").Value = 3934.4
$ws.Range("J132").Value = 336650.84
$ws.Range("K132").Value = 11803.2
$ws.Range("L132").Value = 1009952.52
$ws.Range("M132").Value = -9273.200000000001
$ws.Range("N132").Value = -1015012.52

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 641.53845
$ws.Range("I107").Value = 494.8889
$ws.Range("J107").Value = 971.5
$ws.Range("K107").Value = 1484.6667
$ws.Range("L107").Value = 2914.5
$ws.Range("M107").Value = 435.3333
$ws.Range("N107").Value = -6754.5

$ws.Range("H122").Value = 1377
$ws.Range("I122").Value = 1377
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4131
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -1681

$ws.Range("H123").Value = 40857.125
$ws.Range("I123").Value = 42000
$ws.Range("J123").Value = 39714.25
$ws.Range("K123").Value = 42000
$ws.Range("L123").Value = 39714.25
$ws.Range("M123").Value = -37100
$ws.Range("N123").Value = -49514.25

$ws.Range("H136").Value = 56844.785
$ws.Range("I136").Value = 43285.293
$ws.Range("J136").Value = 81877.69500000001
$ws.Range("K136").Value = 129855.879
$ws.Range("L136").Value = 245633.085
$ws.Range("M136").Value = -127305.879
$ws.Range("N136").Value = -250733.085
